# Updates currentAveragePrice / LevePriceNQ / LevePriceHQ / LeveProfit figures
# (columns H-N) across all eight job sheets, refreshing them with the
# latest Universalis market-board averages pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 417.86667
$ws.Range("I28").Value = 430.14285
$ws.Range("K28").Value = 430.14285
$ws.Range("M28").Value = 54.85714999999999
$ws.Range("H98").Value = 1942099.5
$ws.Range("I98").Value = 1935539.8
$ws.Range("J98").Value = 2003760.6
$ws.Range("K98").Value = 1935539.8
$ws.Range("L98").Value = 2003760.6
$ws.Range("M98").Value = -1934041.8
$ws.Range("N98").Value = -2006756.6
$ws.Range("H101").Value = 941.5454999999999
$ws.Range("I101").Value = 869.625
$ws.Range("K101").Value = 2608.875
$ws.Range("M101").Value = -986.875
$ws.Range("H111").Value = 5129
$ws.Range("I111").Value = 3180.8
$ws.Range("K111").Value = 9542.400000000001
$ws.Range("M111").Value = -6475.400000000001
$ws.Range("H122").Value = 1942099.5
$ws.Range("I122").Value = 1935539.8
$ws.Range("J122").Value = 2003760.6
$ws.Range("K122").Value = 5806619.4
$ws.Range("L122").Value = 6011281.800000001
$ws.Range("M122").Value = -5804169.4
$ws.Range("N122").Value = -6016181.800000001

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3383.125
$ws.Range("I2").Value = 3596.457
$ws.Range("K2").Value = 3596.457
$ws.Range("M2").Value = -3483.457
$ws.Range("H45").Value = 1740.3889
$ws.Range("I45").Value = 1187.4445
$ws.Range("J45").Value = 2293.3333
$ws.Range("K45").Value = 1187.4445
$ws.Range("L45").Value = 2293.3333
$ws.Range("M45").Value = -810.4445000000001
$ws.Range("N45").Value = -3047.3333
$ws.Range("H74").Value = 3432.8357
$ws.Range("I74").Value = 3437.9858
$ws.Range("J74").Value = 3250
$ws.Range("K74").Value = 3437.9858
$ws.Range("L74").Value = 3250
$ws.Range("M74").Value = -2563.9858
$ws.Range("N74").Value = -4998
$ws.Range("H77").Value = 3432.8357
$ws.Range("I77").Value = 3437.9858
$ws.Range("J77").Value = 3250
$ws.Range("K77").Value = 17189.929
$ws.Range("L77").Value = 16250
$ws.Range("M77").Value = -12821.929
$ws.Range("N77").Value = -24986
$ws.Range("H110").Value = 4599.4644
$ws.Range("J110").Value = 10625.875
$ws.Range("L110").Value = 10625.875
$ws.Range("N110").Value = -14715.875
$ws.Range("H116").Value = 3383.125
$ws.Range("I116").Value = 3596.457
$ws.Range("K116").Value = 3596.457
$ws.Range("M116").Value = -1302.457
$ws.Range("H122").Value = 8032.7417
$ws.Range("I122").Value = 5488.148
$ws.Range("K122").Value = 16464.444
$ws.Range("M122").Value = -14014.444
$ws.Range("H132").Value = 2660.6626
$ws.Range("I132").Value = 1891.3867
$ws.Range("K132").Value = 5674.1601
$ws.Range("M132").Value = -3144.1601

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3383.125
$ws.Range("I3").Value = 3596.457
$ws.Range("K3").Value = 3596.457
$ws.Range("M3").Value = -3482.457
$ws.Range("H86").Value = 1567.0358
$ws.Range("I86").Value = 1461.0769
$ws.Range("J86").Value = 2944.5
$ws.Range("K86").Value = 1461.0769
$ws.Range("L86").Value = 2944.5
$ws.Range("M86").Value = -338.0769
$ws.Range("N86").Value = -5190.5
$ws.Range("H89").Value = 1567.0358
$ws.Range("I89").Value = 1461.0769
$ws.Range("J89").Value = 2944.5
$ws.Range("K89").Value = 7305.3845
$ws.Range("L89").Value = 14722.5
$ws.Range("M89").Value = -1689.3845
$ws.Range("N89").Value = -25954.5
$ws.Range("H105").Value = 1617.2759
$ws.Range("I105").Value = 1536.16
$ws.Range("K105").Value = 1536.16
$ws.Range("M105").Value = 210.8399999999999
$ws.Range("H107").Value = 4343.273
$ws.Range("I107").Value = 4296.1113
$ws.Range("K107").Value = 4296.1113
$ws.Range("M107").Value = -2376.1113

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 3883.8333
$ws.Range("J7").Value = 6606.7144
$ws.Range("L7").Value = 6606.7144
$ws.Range("N7").Value = -6832.7144
$ws.Range("H16").Value = 1925.4231
$ws.Range("I16").Value = 1733.4
$ws.Range("K16").Value = 1733.4
$ws.Range("M16").Value = -1446.4
$ws.Range("H22").Value = 44637.582
$ws.Range("I22").Value = 35465.77
$ws.Range("J22").Value = 92331
$ws.Range("K22").Value = 35465.77
$ws.Range("L22").Value = 92331
$ws.Range("M22").Value = -35115.77
$ws.Range("N22").Value = -93031
$ws.Range("H70").Value = 39000
$ws.Range("J70").Value = 39000
$ws.Range("L70").Value = 39000
$ws.Range("N70").Value = -39630
$ws.Range("H73").Value = 39000
$ws.Range("J73").Value = 39000
$ws.Range("L73").Value = 39000
$ws.Range("N73").Value = -41184
$ws.Range("H105").Value = 4243.1113
$ws.Range("I105").Value = 1816.8334
$ws.Range("K105").Value = 1816.8334
$ws.Range("M105").Value = -69.83339999999998
$ws.Range("H107").Value = 4579.8613
$ws.Range("J107").Value = 7497.095
$ws.Range("L107").Value = 7497.095
$ws.Range("N107").Value = -11337.095
$ws.Range("H113").Value = 1925.4231
$ws.Range("I113").Value = 1733.4
$ws.Range("K113").Value = 1733.4
$ws.Range("M113").Value = 436.5999999999999
$ws.Range("H122").Value = 1771.579
$ws.Range("I122").Value = 1532.2354
$ws.Range("J122").Value = 3806
$ws.Range("K122").Value = 4596.706200000001
$ws.Range("L122").Value = 11418
$ws.Range("M122").Value = -2146.706200000001
$ws.Range("N122").Value = -16318
$ws.Range("H132").Value = 3296.8774
$ws.Range("J132").Value = 5211.8
$ws.Range("L132").Value = 15635.4
$ws.Range("N132").Value = -20695.4

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 13889514
$ws.Range("I68").Value = 16667167
$ws.Range("J68").Value = 1249
$ws.Range("K68").Value = 50001501
$ws.Range("L68").Value = 3747
$ws.Range("M68").Value = -50000690
$ws.Range("N68").Value = -5369
$ws.Range("H71").Value = 13889514
$ws.Range("I71").Value = 16667167
$ws.Range("J71").Value = 1249
$ws.Range("K71").Value = 150004503
$ws.Range("L71").Value = 11241
$ws.Range("M71").Value = -150000447
$ws.Range("N71").Value = -19353
$ws.Range("H107").Value = 560.67645
$ws.Range("I107").Value = 589.6667
$ws.Range("J107").Value = 528.0625
$ws.Range("K107").Value = 1769.0001
$ws.Range("L107").Value = 1584.1875
$ws.Range("M107").Value = 150.9999
$ws.Range("N107").Value = -5424.1875

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1748.05
$ws.Range("I122").Value = 1569.2941
$ws.Range("J122").Value = 2761
$ws.Range("K122").Value = 4707.8823
$ws.Range("L122").Value = 8283
$ws.Range("M122").Value = -2257.8823
$ws.Range("N122").Value = -13183
$ws.Range("H136").Value = 36846.375
$ws.Range("J136").Value = 36846.375
$ws.Range("L136").Value = 110539.125
$ws.Range("N136").Value = -115639.125

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1131785.5
$ws.Range("I40").Value = 1271759
$ws.Range("K40").Value = 1271759
$ws.Range("M40").Value = -1271623
$ws.Range("H61").Value = 1108.8823
$ws.Range("I61").Value = 1121.9375
$ws.Range("K61").Value = 1121.9375
$ws.Range("M61").Value = -919.9375
$ws.Range("H113").Value = 1108.8823
$ws.Range("I113").Value = 1121.9375
$ws.Range("K113").Value = 1121.9375
$ws.Range("M113").Value = 1048.0625
$ws.Range("H122").Value = 3091.9524
$ws.Range("J122").Value = 5268.3335
$ws.Range("L122").Value = 15805.0005
$ws.Range("N122").Value = -20705.0005

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 14488.2
$ws.Range("I14").Value = 20072
$ws.Range("K14").Value = 20072
$ws.Range("M14").Value = -19904
$ws.Range("H107").Value = 709.7778
$ws.Range("I107").Value = 709.75
$ws.Range("J107").Value = 709.8
$ws.Range("K107").Value = 2129.25
$ws.Range("L107").Value = 2129.4
$ws.Range("M107").Value = -209.25
$ws.Range("N107").Value = -5969.4
$ws.Range("H122").Value = 1548.9688
$ws.Range("I122").Value = 1323.2
$ws.Range("K122").Value = 3969.6
$ws.Range("M122").Value = -1519.6
$ws.Range("H132").Value = 1428.91
$ws.Range("I132").Value = 755.6393399999999
$ws.Range("J132").Value = 2481.9744
$ws.Range("K132").Value = 2266.91802
$ws.Range("L132").Value = 7445.9232
$ws.Range("N132").Value = -12505.9232

Write-Host "Updated market-price figures across all sheets."
